# Update the three result sheets (sigma_010, sigma_025, sigma_050) with new
# NLM-LBP run values: the "Rows" index column now starts at 1 (was 0), and
# the Noisy/NLM-LBP measurement columns are refreshed with new results.

$wb = $excel.ActiveWorkbook

$sheetNames = @("sigma_010", "sigma_025", "sigma_050")

$data = @{
    "sigma_010" = @(
        @(1,  27.8187978278512,  28.57271655009389),
        @(2,  27.84397970724693, 28.57640606645103),
        @(3,  27.80539676323497, 28.54546697135685),
        @(4,  27.82297592306,    28.56299240347486),
        @(5,  27.8354225098857,  28.58581980844269),
        @(6,  27.81040955362748, 28.56310078290642),
        @(7,  27.86156351606635, 28.56819141537427),
        @(8,  27.80022297296503, 28.55911095503167),
        @(9,  27.85644969396837, 28.57217316597544),
        @(10, 27.82439059371321, 28.58951576450293)
    )
    "sigma_025" = @(
        @(1,  19.72936212069984, 25.08057875576624),
        @(2,  19.74855385633919, 25.0918863595157),
        @(3,  19.72051751290509, 25.08626257521509),
        @(4,  19.74379053495597, 25.05167107214245),
        @(5,  19.74566073564925, 25.07339840203907),
        @(6,  19.73140618631875, 25.08198151545008),
        @(7,  19.72768394373398, 25.09057151668794),
        @(8,  19.72374773196849, 25.07413790502166),
        @(9,  19.71064740330682, 25.05608468275296),
        @(10, 19.75112047036913, 25.09272692445177)
    )
    "sigma_050" = @(
        @(1,  14.67819158410424, 20.14182173884474),
        @(2,  14.67184926604159, 20.15149388658961),
        @(3,  14.66405295758777, 20.16994143319047),
        @(4,  14.67252334041433, 20.19408792556247),
        @(5,  14.6783063978622,  20.1678816525956),
        @(6,  14.66222145325498, 20.17787938655982),
        @(7,  14.67102942094896, 20.15681267141502),
        @(8,  14.68480209579522, 20.14666221504719),
        @(9,  14.67869558954255, 20.15284212539778),
        @(10, 14.69907134584581, 20.17954868207924)
    )
}

$means = @{
    "sigma_010" = @(27.82796090616192, 28.569549388361)
    "sigma_025" = @(19.73324904962465, 25.0779299709043)
    "sigma_050" = @(14.67607434513976, 20.16389717172819)
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    for ($i = 0; $i -lt $data[$name].Count; $i++) {
        $row = 2 + $i
        $vals = $data[$name][$i]
        $ws.Cells.Item($row, 1).Value = $vals[0]
        $ws.Cells.Item($row, 2).Value = $vals[1]
        $ws.Cells.Item($row, 3).Value = $vals[2]
    }

    $meanVals = $means[$name]
    $ws.Cells.Item(12, 2).Value = $meanVals[0]
    $ws.Cells.Item(12, 3).Value = $meanVals[1]
}
